$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Row 2 (C40) gets the old Row 3 content for B, C, D, F
$ws.Cells.Item(2,2).Value = 'Media and Communications Senior Manager, GCoM'
$ws.Cells.Item(2,3).Value = 'Global Covenant of Mayors - Brazil, South Africa'
$ws.Cells.Item(2,4).Value = 'N/A'
$ws.Cells.Item(2,6).Value = 'https://c40.bamboohr.com/careers/705'

# Step 2: Insert 3 new rows before row 4 (pushes old row4.. down by 3)
$ws.Rows("4:6").Insert()
$ws.Rows("4:6").RowHeight = 80

# Step 3: Overwrite row 3 in place with the new TA-6822 entry; clear C3
$ws.Cells.Item(3,1).Value = 'DevelopmentAid'
$ws.Cells.Item(3,2).Value = 'TA-6822 IND: Support for Strengthening Multimodal and Integrated Logistics Ecosystem - Development of Sectoral Plan for Efficient Logistics (SPEL) for Wheat, Rice and Millets under Public Distribution System - Agri-market Assessment Expert (55154-002)'
$ws.Cells.Item(3,3).Value = ''
$ws.Cells.Item(3,4).Value = 'Governance, Learning'
$ws.Cells.Item(3,6).Value = 'https://www.developmentaid.org/tenders/view/1607893/ta-6822-ind-support-for-strengthening-multimodal-and-integrated-logistics-ecosystem-development-of-s'

# Step 4: Fill newly inserted rows 4, 5, 6
$ws.Cells.Item(4,1).Value = 'DevelopmentAid'
$ws.Cells.Item(4,2).Value = 'RFP - Final Evaluation of the Bhoomi Ka Programme'
$ws.Cells.Item(4,4).Value = 'Governance'
$ws.Cells.Item(4,6).Value = 'https://www.developmentaid.org/tenders/view/1612729/rfp-final-evaluation-of-the-bhoomi-ka-programme'

$ws.Cells.Item(5,1).Value = 'DevelopmentAid'
$ws.Cells.Item(5,2).Value = 'RFP- for Hiring a Resource Person/Agency Development of MEAL Questionnaires and Processes for CommCare Tool'
$ws.Cells.Item(5,4).Value = 'Governance, Learning'
$ws.Cells.Item(5,6).Value = 'https://www.developmentaid.org/tenders/view/1612750/rfp-for-hiring-a-resource-personagency-development-of-meal-questionnaires-and-processes-for-commcare'

$ws.Cells.Item(6,1).Value = 'DevelopmentAid'
$ws.Cells.Item(6,2).Value = 'P178254- Kerala Climate Resilient Agri- Value Chain Modernization (KERA) Project - Procurement Plan (State Project Management Unit (Kerala))'
$ws.Cells.Item(6,4).Value = 'Governance, Climate'
$ws.Cells.Item(6,6).Value = 'https://www.developmentaid.org/tenders/view/1612677/p178254-kerala-climate-resilient-agri-value-chain-modernization-kera-project-procurement-plan-state'

# Step 5: Delete the two rows pushed off the bottom (old rows 61 and 62 now at 64,65)
$ws.Rows("64:65").Delete()